$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.934.49"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "3.392.05"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.30"
$ws.Range("E5").Value = "  -1.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.60"
$ws.Range("E6").Value = "  -2.22%  "
$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").Value = "3.391.26"
$ws.Range("E7").Value = "  -1.65%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.475"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.51"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("E12").Value = "  +1.79%  "
$ws.Range("D13").Value = "3.970.89"
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.48"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("D17").Value = "3.389.25"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").Value = "60.979.23"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.22"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.98"
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.00"
$ws.Range("E21").Value = "  -5.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "384.63"
$ws.Range("E22").Value = "  -2.60%  "
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.05"
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000117"
$ws.Range("E26").Value = "  -5.17%  "
$ws.Range("D27").Value = "3.526.14"
$ws.Range("E27").Value = "  -1.82%  "
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -2.39%  "
$ws.Range("E31").Value = "  -2.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.15"
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.67"
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.01"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "165.72"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("D38").Value = "3.423.36"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.99"
$ws.Range("E39").Value = "  -2.57%  "
$ws.Range("E40").Value = "  -4.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "28.06"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0775"
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.780"
$ws.Range("E44").Value = "  -2.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.98"
$ws.Range("E45").Value = "  -0.65%  "
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.67"
$ws.Range("E47").Value = "  -3.62%  "
$ws.Range("E48").Value = "  -2.49%  "
$ws.Range("D49").Value = "2.488.41"
$ws.Range("E49").Value = "  -4.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.43"
$ws.Range("E50").Value = "  +1.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.82"
$ws.Range("E51").Value = "  -1.21%  "
